$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2-79).
# Bump the serial value by one day (45188 -> 45189) while preserving the
# existing numeric type and cell style/format.
for ($row = 2; $row -le 79; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
